# Bugfix for the naive forecaster component: refresh the YoY forecast vectors
# (dates/years/diffs shift by one period and are recomputed; final stale row dropped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 39583
$ws.Cells.Item(2, 2).Value = 2008
$ws.Cells.Item(2, 4).Value = 2009
$ws.Cells.Item(2, 5).Value = 0.07975464681371225
# Row 3
$ws.Cells.Item(3, 1).Value = 39765
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = -4.700509864312973
# Row 4
$ws.Cells.Item(4, 1).Value = 39948
$ws.Cells.Item(4, 2).Value = 2009
$ws.Cells.Item(4, 3).Value = -3.017636378339217
$ws.Cells.Item(4, 4).Value = 2010
$ws.Cells.Item(4, 5).Value = -3.229247082222797
# Row 5
$ws.Cells.Item(5, 1).Value = 40130
$ws.Cells.Item(5, 2).Value = 2009
$ws.Cells.Item(5, 3).Value = -1.31761954538927
$ws.Cells.Item(5, 4).Value = 2010
$ws.Cells.Item(5, 5).Value = -0.01655958389530365
# Row 6
$ws.Cells.Item(6, 1).Value = 40310
$ws.Cells.Item(6, 2).Value = 2010
$ws.Cells.Item(6, 3).Value = 0.003352386816724007
$ws.Cells.Item(6, 4).Value = 2011
$ws.Cells.Item(6, 5).Value = -1.461031976610316
# Row 7
$ws.Cells.Item(7, 1).Value = 40494
$ws.Cells.Item(7, 2).Value = 2010
$ws.Cells.Item(7, 3).Value = 2.771597318554297
$ws.Cells.Item(7, 4).Value = 2011
$ws.Cells.Item(7, 5).Value = 3.579142225970444
# Row 8
$ws.Cells.Item(8, 1).Value = 40676
$ws.Cells.Item(8, 2).Value = 2011
$ws.Cells.Item(8, 3).Value = 3.925837669383347
$ws.Cells.Item(8, 4).Value = 2012
$ws.Cells.Item(8, 5).Value = 2.641604203902781
# Row 9
$ws.Cells.Item(9, 1).Value = 40862
$ws.Cells.Item(9, 2).Value = 2011
$ws.Cells.Item(9, 3).Value = 1.799362536952542
$ws.Cells.Item(9, 4).Value = 2012
$ws.Cells.Item(9, 5).Value = -0.289184878867832
# Row 10
$ws.Cells.Item(10, 1).Value = 41044
$ws.Cells.Item(10, 2).Value = 2012
$ws.Cells.Item(10, 3).Value = 0.2381541440396262
$ws.Cells.Item(10, 4).Value = 2013
$ws.Cells.Item(10, 5).Value = 1.60268309892857
# Row 11
$ws.Cells.Item(11, 1).Value = 41228
$ws.Cells.Item(11, 2).Value = 2012
$ws.Cells.Item(11, 3).Value = 2.123182427147152
$ws.Cells.Item(11, 4).Value = 2013
$ws.Cells.Item(11, 5).Value = 5.963492031746176
# Row 12
$ws.Cells.Item(12, 1).Value = 41409
$ws.Cells.Item(12, 2).Value = 2013
$ws.Cells.Item(12, 3).Value = 4.993892964711621
$ws.Cells.Item(12, 4).Value = 2014
$ws.Cells.Item(12, 5).Value = 2.260118192030736
# Row 13
$ws.Cells.Item(13, 1).Value = 41592
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = 7.317297369612819
$ws.Cells.Item(13, 4).Value = 2014
$ws.Cells.Item(13, 5).Value = 7.523777575896196
# Row 14
$ws.Cells.Item(14, 1).Value = 41774
$ws.Cells.Item(14, 2).Value = 2014
$ws.Cells.Item(14, 3).Value = 6.711795724673664
$ws.Cells.Item(14, 4).Value = 2015
$ws.Cells.Item(14, 5).Value = 6.409878804372982
# Row 15
$ws.Cells.Item(15, 1).Value = 41957
$ws.Cells.Item(15, 2).Value = 2014
$ws.Cells.Item(15, 3).Value = 4.260319658857736
$ws.Cells.Item(15, 4).Value = 2015
$ws.Cells.Item(15, 5).Value = 2.532215190177589
# Row 16
$ws.Cells.Item(16, 1).Value = 42137
$ws.Cells.Item(16, 2).Value = 2015
$ws.Cells.Item(16, 3).Value = 0.5121603413743347
$ws.Cells.Item(16, 4).Value = 2016
$ws.Cells.Item(16, 5).Value = 3.290935868252554
# Row 17
$ws.Cells.Item(17, 1).Value = 42321
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = -0.05262415810141086
$ws.Cells.Item(17, 4).Value = 2016
$ws.Cells.Item(17, 5).Value = 2.051185924063259
# Row 18
$ws.Cells.Item(18, 1).Value = 42503
$ws.Cells.Item(18, 2).Value = 2016
$ws.Cells.Item(18, 3).Value = 1.745565778643887
$ws.Cells.Item(18, 4).Value = 2017
$ws.Cells.Item(18, 5).Value = 0.7985845180024986
# Row 19
$ws.Cells.Item(19, 1).Value = 42689
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = 1.459778471779982
$ws.Cells.Item(19, 4).Value = 2017
$ws.Cells.Item(19, 5).Value = 0.4575538530338541
# Row 20
$ws.Cells.Item(20, 1).Value = 42867
$ws.Cells.Item(20, 2).Value = 2017
$ws.Cells.Item(20, 3).Value = 2.687500891103922
$ws.Cells.Item(20, 4).Value = 2018
$ws.Cells.Item(20, 5).Value = 1.922191950024699
# Row 21
$ws.Cells.Item(21, 1).Value = 43053
$ws.Cells.Item(21, 2).Value = 2017
$ws.Cells.Item(21, 3).Value = 3.002208343813528
$ws.Cells.Item(21, 4).Value = 2018
$ws.Cells.Item(21, 5).Value = 2.600569166164624
# Row 22
$ws.Cells.Item(22, 1).Value = 43145
$ws.Cells.Item(22, 2).Value = 2018
$ws.Cells.Item(22, 3).Value = 3.339205815020496
$ws.Cells.Item(22, 4).Value = 2019
$ws.Cells.Item(22, 5).Value = 2.646788941483735
# Row 23
$ws.Cells.Item(23, 1).Value = 43235
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = 3.654655474034474
$ws.Cells.Item(23, 4).Value = 2019
$ws.Cells.Item(23, 5).Value = 3.068403604789749
# Row 24
$ws.Cells.Item(24, 1).Value = 43326
$ws.Cells.Item(24, 2).Value = 2018
$ws.Cells.Item(24, 3).Value = 4.036117574265741
$ws.Cells.Item(24, 4).Value = 2019
$ws.Cells.Item(24, 5).Value = 3.628675245064317
# Row 25
$ws.Cells.Item(25, 1).Value = 43418
$ws.Cells.Item(25, 2).Value = 2018
$ws.Cells.Item(25, 3).Value = 4.020433260014977
$ws.Cells.Item(25, 4).Value = 2019
$ws.Cells.Item(25, 5).Value = 3.605726003451304
# Row 26
$ws.Cells.Item(26, 1).Value = 43510
$ws.Cells.Item(26, 2).Value = 2019
$ws.Cells.Item(26, 3).Value = 3.813466308501412
$ws.Cells.Item(26, 4).Value = 2020
$ws.Cells.Item(26, 5).Value = 3.884173085820986
# Row 27
$ws.Cells.Item(27, 1).Value = 43600
$ws.Cells.Item(27, 2).Value = 2019
$ws.Cells.Item(27, 3).Value = 3.712036718632117
$ws.Cells.Item(27, 4).Value = 2020
$ws.Cells.Item(27, 5).Value = 3.908921577463587
# Row 28
$ws.Cells.Item(28, 1).Value = 43691
$ws.Cells.Item(28, 2).Value = 2019
$ws.Cells.Item(28, 3).Value = 3.551357200054261
$ws.Cells.Item(28, 4).Value = 2020
$ws.Cells.Item(28, 5).Value = 3.524277826276134
# Row 29
$ws.Cells.Item(29, 1).Value = 43783
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = 3.53224976671227
$ws.Cells.Item(29, 4).Value = 2020
$ws.Cells.Item(29, 5).Value = 3.490656491795074
# Row 30
$ws.Cells.Item(30, 1).Value = 43875
$ws.Cells.Item(30, 2).Value = 2020
$ws.Cells.Item(30, 3).Value = 2.732790977059629
$ws.Cells.Item(30, 4).Value = 2021
$ws.Cells.Item(30, 5).Value = 3.412957258051663
# Row 31
$ws.Cells.Item(31, 1).Value = 43966
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = 2.849400388885992
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = 3.535456592693387
# Row 32
$ws.Cells.Item(32, 1).Value = 44068
$ws.Cells.Item(32, 2).Value = 2020
$ws.Cells.Item(32, 3).Value = -0.985458715495402
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = -3.099271113627677
# Row 33
$ws.Cells.Item(33, 1).Value = 44159
$ws.Cells.Item(33, 2).Value = 2020
$ws.Cells.Item(33, 3).Value = -0.985458715495402
$ws.Cells.Item(33, 4).Value = 2021
$ws.Cells.Item(33, 5).Value = -2.347097924577757
# Row 34
$ws.Cells.Item(34, 1).Value = 44251
$ws.Cells.Item(34, 2).Value = 2021
$ws.Cells.Item(34, 3).Value = -5.285337128797329
$ws.Cells.Item(34, 4).Value = 2022
$ws.Cells.Item(34, 5).Value = -2.69436300797079
# Row 35
$ws.Cells.Item(35, 1).Value = 44341
$ws.Cells.Item(35, 2).Value = 2021
$ws.Cells.Item(35, 3).Value = -4.741003096464214
$ws.Cells.Item(35, 4).Value = 2022
$ws.Cells.Item(35, 5).Value = -2.156795995006056
# Row 36
$ws.Cells.Item(36, 1).Value = 44432
$ws.Cells.Item(36, 2).Value = 2021
$ws.Cells.Item(36, 3).Value = -4.365687260408224
$ws.Cells.Item(36, 4).Value = 2022
$ws.Cells.Item(36, 5).Value = -0.9191921099315992
# Row 37
$ws.Cells.Item(37, 1).Value = 44525
$ws.Cells.Item(37, 2).Value = 2021
$ws.Cells.Item(37, 3).Value = -4.365687260408224
$ws.Cells.Item(37, 4).Value = 2022
$ws.Cells.Item(37, 5).Value = -0.1803381976702711
# Row 38
$ws.Cells.Item(38, 1).Value = 44617
$ws.Cells.Item(38, 2).Value = 2022
$ws.Cells.Item(38, 3).Value = 1.871837441670499
$ws.Cells.Item(38, 4).Value = 2023
$ws.Cells.Item(38, 5).Value = -0.8573220257725445
# Row 39
$ws.Cells.Item(39, 1).Value = 44706
$ws.Cells.Item(39, 2).Value = 2022
$ws.Cells.Item(39, 3).Value = 1.194925448553708
$ws.Cells.Item(39, 4).Value = 2023
$ws.Cells.Item(39, 5).Value = -2.068675356622807
# Row 40
$ws.Cells.Item(40, 1).Value = 44798
$ws.Cells.Item(40, 2).Value = 2022
$ws.Cells.Item(40, 3).Value = 1.27347919322387
$ws.Cells.Item(40, 4).Value = 2023
$ws.Cells.Item(40, 5).Value = -2.137023292796481
# Row 41
$ws.Cells.Item(41, 1).Value = 44890
$ws.Cells.Item(41, 2).Value = 2022
$ws.Cells.Item(41, 3).Value = 1.27347919322387
$ws.Cells.Item(41, 4).Value = 2023
$ws.Cells.Item(41, 5).Value = -1.152671696465724
# Row 42
$ws.Cells.Item(42, 1).Value = 44981
$ws.Cells.Item(42, 2).Value = 2023
$ws.Cells.Item(42, 3).Value = -0.507642354784088
$ws.Cells.Item(42, 4).Value = 2024
$ws.Cells.Item(42, 5).Value = 0.07646803116447831
# Row 43
$ws.Cells.Item(43, 1).Value = 45071
$ws.Cells.Item(43, 2).Value = 2023
$ws.Cells.Item(43, 3).Value = -1.084365158506884
$ws.Cells.Item(43, 4).Value = 2024
$ws.Cells.Item(43, 5).Value = -1.089896342664354
# Row 44
$ws.Cells.Item(44, 1).Value = 45163
$ws.Cells.Item(44, 2).Value = 2023
$ws.Cells.Item(44, 3).Value = -1.339436245206127
$ws.Cells.Item(44, 4).Value = 2024
$ws.Cells.Item(44, 5).Value = -1.277727682704721
# Row 45
$ws.Cells.Item(45, 1).Value = 45254
$ws.Cells.Item(45, 2).Value = 2023
$ws.Cells.Item(45, 3).Value = -1.339436245206127
$ws.Cells.Item(45, 4).Value = 2024
$ws.Cells.Item(45, 5).Value = -2.785556326028149
# Row 46
$ws.Cells.Item(46, 1).Value = 45345
$ws.Cells.Item(46, 2).Value = 2024
$ws.Cells.Item(46, 3).Value = -2.87408779878463
$ws.Cells.Item(46, 4).Value = 2025
$ws.Cells.Item(46, 5).Value = -1.155307395925487
# Row 47
$ws.Cells.Item(47, 1).Value = 45436
$ws.Cells.Item(47, 2).Value = 2024
$ws.Cells.Item(47, 3).Value = -3.40787540386569
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = -2.191935020614488
# Row 48
$ws.Cells.Item(48, 1).Value = 45534
$ws.Cells.Item(48, 2).Value = 2024
$ws.Cells.Item(48, 3).Value = -3.451527003230626
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = -2.998472503487815
# Row 49
$ws.Cells.Item(49, 1).Value = 45618
$ws.Cells.Item(49, 2).Value = 2024
$ws.Cells.Item(49, 3).Value = -3.451527003230626
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = -2.452009576682213
# Row 50
$ws.Cells.Item(50, 1).Value = 45713
$ws.Cells.Item(50, 2).Value = 2025
$ws.Cells.Item(50, 3).Value = -2.164008261629446
$ws.Cells.Item(50, 4).Value = 2026
$ws.Cells.Item(50, 5).Value = -2.337170009804157
# Row 51
$ws.Cells.Item(51, 1).Value = 45800
$ws.Cells.Item(51, 2).Value = 2025
$ws.Cells.Item(51, 3).Value = -1.853660925652212
$ws.Cells.Item(51, 4).Value = 2026
$ws.Cells.Item(51, 5).Value = -1.878672029998096
# Row 52
$ws.Cells.Item(52, 1).Value = 45891
$ws.Cells.Item(52, 2).Value = 2025
$ws.Cells.Item(52, 3).Value = -1.75044229618867
$ws.Cells.Item(52, 4).Value = 2026
$ws.Cells.Item(52, 5).Value = -1.760724207457021

# Drop the now-stale trailing row (data window shifted forward by one period)
$ws.Rows.Item(53).Delete()
